# Applies numeric cell updates to the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets
# as captured by the source diff (scheduled-runner price/profit refresh).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 2
$ws.Range("H2").Value = 95.666664
$ws.Range("I2").Value = 84.8
$ws.Range("J2").Value = 150
$ws.Range("K2").Value = 84.8
$ws.Range("L2").Value = 150
$ws.Range("M2").Value = 28.2
$ws.Range("N2").Value = -376

# Row 4
$ws.Range("H4").Value = 419.2
$ws.Range("I4").Value = 419.2
$ws.Range("K4").Value = 419.2
$ws.Range("M4").Value = -305.2

# Row 5
$ws.Range("H5").Value = 82.375
$ws.Range("I5").Value = 82.375
$ws.Range("K5").Value = 82.375
$ws.Range("M5").Value = 32.625

# Row 9
$ws.Range("H9").Value = 281.2857
$ws.Range("I9").Value = 37.5
$ws.Range("K9").Value = 37.5
$ws.Range("M9").Value = 131.5

# Row 19
$ws.Range("H19").Value = 1075.4445
$ws.Range("I19").Value = 955
$ws.Range("J19").Value = 1316.3334
$ws.Range("K19").Value = 955
$ws.Range("L19").Value = 1316.3334
$ws.Range("M19").Value = -780
$ws.Range("N19").Value = -1666.3334

# Row 40
$ws.Range("H40").Value = 61527.598
$ws.Range("J40").Value = 3947.617
$ws.Range("L40").Value = 3947.617
$ws.Range("N40").Value = -4297.617

# Row 132
$ws.Range("H132").Value = 5413.6787
$ws.Range("I132").Value = 5703.32
$ws.Range("K132").Value = 17109.96
$ws.Range("M132").Value = -14579.96

# Row 137
$ws.Range("H137").Value = 51534.41
$ws.Range("I137").Value = 84219.7
$ws.Range("J137").Value = 4841.143
$ws.Range("K137").Value = 252659.1
$ws.Range("L137").Value = 14523.429
$ws.Range("M137").Value = -250109.1
$ws.Range("N137").Value = -19623.429

$ws = $wb.Worksheets.Item("ARM")
# Row 13
$ws.Range("H13").Value = 1669366.6
$ws.Range("J13").Value = 4050
$ws.Range("L13").Value = 4050
$ws.Range("N13").Value = -4338

# Row 32
$ws.Range("H32").Value = 1957.1
$ws.Range("I32").Value = 1792.6735
$ws.Range("J32").Value = 10014
$ws.Range("K32").Value = 1792.6735
$ws.Range("L32").Value = 10014
$ws.Range("M32").Value = -1505.6735
$ws.Range("N32").Value = -10588

# Row 63
$ws.Range("H63").Value = 120042
$ws.Range("I63").Value = 2242.3333
$ws.Range("J63").Value = 237841.67
$ws.Range("K63").Value = 2242.3333
$ws.Range("L63").Value = 237841.67
$ws.Range("M63").Value = -1556.3333
$ws.Range("N63").Value = -239213.67

# Row 66
$ws.Range("H66").Value = 120042
$ws.Range("I66").Value = 2242.3333
$ws.Range("J66").Value = 237841.67
$ws.Range("K66").Value = 11211.6665
$ws.Range("L66").Value = 1189208.35
$ws.Range("M66").Value = -7779.666499999999
$ws.Range("N66").Value = -1196072.35

# Row 97
$ws.Range("H97").Value = 1256.9395
$ws.Range("I97").Value = 1127.5385
$ws.Range("J97").Value = 1737.5714
$ws.Range("K97").Value = 1127.5385
$ws.Range("L97").Value = 1737.5714
$ws.Range("M97").Value = -631.5385000000001
$ws.Range("N97").Value = -2729.5714

# Row 122
$ws.Range("H122").Value = 3923.081
$ws.Range("I122").Value = 2431.7036
$ws.Range("J122").Value = 7949.8
$ws.Range("K122").Value = 7295.110799999999
$ws.Range("L122").Value = 23849.4
$ws.Range("M122").Value = -4845.110799999999
$ws.Range("N122").Value = -28749.4

# Row 132
$ws.Range("H132").Value = 2644.682
$ws.Range("I132").Value = 2552.0571
$ws.Range("J132").Value = 3004.889
$ws.Range("K132").Value = 7656.1713
$ws.Range("L132").Value = 9014.667000000001
$ws.Range("M132").Value = -5126.1713
$ws.Range("N132").Value = -14074.667

$ws = $wb.Worksheets.Item("BSM")
# Row 22
$ws.Range("H22").Value = 331.66666
$ws.Range("I22").Value = 331.66666
$ws.Range("K22").Value = 331.66666
$ws.Range("M22").Value = -158.66666

# Row 58
$ws.Range("H58").Value = 8580
$ws.Range("I58").Value = 0
$ws.Range("J58").Value = 8580
$ws.Range("K58").Value = 0
$ws.Range("L58").Value = 8580
$ws.Range("M58").ClearContents()
$ws.Range("N58").Value = -9168

# Row 60
$ws.Range("H60").Value = 105052
$ws.Range("J60").Value = 105052
$ws.Range("L60").Value = 105052
$ws.Range("N60").Value = -106250

# Row 87
$ws.Range("H87").Value = 149000
$ws.Range("I87").Value = 0
$ws.Range("K87").Value = 0
$ws.Range("M87").ClearContents()

# Row 90
$ws.Range("H90").Value = 149000
$ws.Range("I90").Value = 0
$ws.Range("K90").Value = 0
$ws.Range("M90").ClearContents()

# Row 117
$ws.Range("H117").Value = 117990
$ws.Range("I117").Value = 0
$ws.Range("J117").Value = 117990
$ws.Range("K117").Value = 0
$ws.Range("L117").Value = 117990
$ws.Range("M117").ClearContents()
$ws.Range("N117").Value = -127168

$ws = $wb.Worksheets.Item("CRP")
# Row 7
$ws.Range("H7").Value = 60
$ws.Range("I7").Value = 64
$ws.Range("K7").Value = 64
$ws.Range("M7").Value = 49

# Row 122
$ws.Range("H122").Value = 1320.6666
$ws.Range("I122").Value = 1344.8
$ws.Range("J122").Value = 1200
$ws.Range("K122").Value = 4034.4
$ws.Range("L122").Value = 3600
$ws.Range("M122").Value = -1584.4
$ws.Range("N122").Value = -8500

# Row 134
$ws.Range("H134").Value = 3036.1428
$ws.Range("I134").Value = 2388.762
$ws.Range("K134").Value = 7166.286
$ws.Range("M134").Value = -4631.286

$ws = $wb.Worksheets.Item("CUL")
# Row 4
$ws.Range("H4").Value = 77172200
$ws.Range("I4").Value = 92141650
$ws.Range("J4").Value = 63999084
$ws.Range("K4").Value = 276424950
$ws.Range("L4").Value = 191997252
$ws.Range("M4").Value = -276424838
$ws.Range("N4").Value = -191997476

# Row 11
$ws.Range("H11").Value = 672516.3
$ws.Range("I11").Value = 842045.5
$ws.Range("K11").Value = 2526136.5
$ws.Range("M11").Value = -2525996.5

# Row 17
$ws.Range("H17").Value = 2131.4285
$ws.Range("I17").Value = 3433.3333
$ws.Range("J17").Value = 1155
$ws.Range("K17").Value = 10299.9999
$ws.Range("L17").Value = 3465
$ws.Range("M17").Value = -10130.9999
$ws.Range("N17").Value = -3803

# Row 24
$ws.Range("H24").Value = 5465
$ws.Range("I24").Value = 699.5
$ws.Range("K24").Value = 2098.5
$ws.Range("M24").Value = -1868.5

# Row 34
$ws.Range("H34").Value = 696.625
$ws.Range("I34").Value = 103.1
$ws.Range("J34").Value = 1685.8334
$ws.Range("K34").Value = 309.3
$ws.Range("L34").Value = 5057.5002
$ws.Range("M34").Value = -225.3
$ws.Range("N34").Value = -5225.5002

# Row 35
$ws.Range("H35").Value = 2002
$ws.Range("I35").Value = 2002
$ws.Range("J35").Value = 0
$ws.Range("K35").Value = 6006
$ws.Range("L35").Value = 0
$ws.Range("M35").Value = -5718
$ws.Range("N35").ClearContents()

# Row 46
$ws.Range("H46").Value = 0
$ws.Range("J46").Value = 0
$ws.Range("L46").Value = 0
$ws.Range("N46").ClearContents()

# Row 51
$ws.Range("H51").Value = 1165.7
$ws.Range("I51").Value = 368.66666
$ws.Range("J51").Value = 1507.2858
$ws.Range("K51").Value = 1105.99998
$ws.Range("L51").Value = 4521.857400000001
$ws.Range("M51").Value = -645.9999800000001
$ws.Range("N51").Value = -5441.857400000001

# Row 57
$ws.Range("H57").Value = 5401.2
$ws.Range("I57").Value = 3666.6667
$ws.Range("J57").Value = 5834.8335
$ws.Range("K57").Value = 11000.0001
$ws.Range("L57").Value = 17504.5005
$ws.Range("M57").Value = -10441.0001
$ws.Range("N57").Value = -18622.5005

# Row 60
$ws.Range("H60").Value = 407.66666
$ws.Range("I60").Value = 581.6667
$ws.Range("J60").Value = 233.66667
$ws.Range("K60").Value = 1745.0001
$ws.Range("L60").Value = 701.00001
$ws.Range("M60").Value = -1494.0001
$ws.Range("N60").Value = -1203.00001

# Row 113
$ws.Range("H113").Value = 1160.7916
$ws.Range("J113").Value = 1203.1818
$ws.Range("L113").Value = 3609.5454
$ws.Range("N113").Value = -7949.5454

# Row 117
$ws.Range("H117").Value = 2130.72
$ws.Range("J117").Value = 2215.6086
$ws.Range("L117").Value = 6646.825800000001
$ws.Range("N117").Value = -13530.8258

# Row 132
$ws.Range("H132").Value = 467.6
$ws.Range("I132").Value = 467.6
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 4208.400000000001
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -1678.400000000001
$ws.Range("N132").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
# Row 102
$ws.Range("H102").Value = 2986
$ws.Range("I102").Value = 2986
$ws.Range("K102").Value = 2986
$ws.Range("M102").Value = -1364

# Row 132
$ws.Range("H132").Value = 4379
$ws.Range("I132").Value = 4928.2856
$ws.Range("J132").Value = 2456.5
$ws.Range("K132").Value = 14784.8568
$ws.Range("L132").Value = 7369.5
$ws.Range("M132").Value = -12254.8568
$ws.Range("N132").Value = -12429.5

# Row 134
$ws.Range("H134").Value = 59748.75
$ws.Range("J134").Value = 59748.75
$ws.Range("L134").Value = 179246.25
$ws.Range("N134").Value = -184316.25

$ws = $wb.Worksheets.Item("LTW")
# Row 16
$ws.Range("H16").Value = 2199.5557
$ws.Range("J16").Value = 0
$ws.Range("L16").Value = 0
$ws.Range("N16").ClearContents()

# Row 22
$ws.Range("H22").Value = 3488.9312
$ws.Range("I22").Value = 2559.6924
$ws.Range("J22").Value = 4243.9375
$ws.Range("K22").Value = 2559.6924
$ws.Range("L22").Value = 4243.9375
$ws.Range("M22").Value = -2264.6924
$ws.Range("N22").Value = -4833.9375

# Row 27
$ws.Range("H27").Value = 3488.9312
$ws.Range("I27").Value = 2559.6924
$ws.Range("J27").Value = 4243.9375
$ws.Range("K27").Value = 2559.6924
$ws.Range("L27").Value = 4243.9375
$ws.Range("M27").Value = -2452.6924
$ws.Range("N27").Value = -4457.9375

# Row 40
$ws.Range("H40").Value = 6149
$ws.Range("I40").Value = 5507.3335
$ws.Range("K40").Value = 5507.3335
$ws.Range("M40").Value = -5371.3335

# Row 46
$ws.Range("H46").Value = 2470.8
$ws.Range("I46").Value = 670
$ws.Range("J46").Value = 3071.0667
$ws.Range("K46").Value = 670
$ws.Range("L46").Value = 3071.0667
$ws.Range("M46").Value = -482
$ws.Range("N46").Value = -3447.0667

# Row 100
$ws.Range("H100").Value = 3422.0667
$ws.Range("I100").Value = 3304
$ws.Range("J100").Value = 3599.1667
$ws.Range("K100").Value = 3304
$ws.Range("L100").Value = 3599.1667
$ws.Range("M100").Value = -2763
$ws.Range("N100").Value = -4681.1667

$ws = $wb.Worksheets.Item("WVR")
# Row 62
$ws.Range("H62").Value = 5100.5
$ws.Range("I62").Value = 4100.4
$ws.Range("K62").Value = 4100.4
$ws.Range("M62").Value = -3476.4

# Row 65
$ws.Range("H65").Value = 5100.5
$ws.Range("I65").Value = 4100.4
$ws.Range("K65").Value = 20502
$ws.Range("M65").Value = -17382

# Row 122
$ws.Range("H122").Value = 3771.074
$ws.Range("I122").Value = 1643.55
$ws.Range("K122").Value = 4930.65
$ws.Range("M122").Value = -2480.65

